# Insert two new rows at 806-807 (pushes existing rows 806..899 down to 808..901)
# and populate them with the new weekly price data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("806:807").Insert()

# New row 806
$ws.Range("A806").Value = 4
$ws.Range("B806").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C806").Value = "Los Lagos"
$ws.Range("D806").Value = 44918
$ws.Range("E806").Value = 10
$ws.Range("F806").Value = 100112020
$ws.Range("G806").Value = "Tomate"
$ws.Range("H806").Value = "Larga vida"
$ws.Range("I806").Value = "Extra"
$ws.Range("J806").Value = 250
$ws.Range("K806").Value = 24000
$ws.Range("L806").Value = 24000
$ws.Range("M806").Value = 24000
$ws.Range("N806").Value = "$/bandeja 18 kilos"
$ws.Range("O806").Value = "Provincia de Quillota"
$ws.Range("P806").Value = 1333
$ws.Range("Q806").Value = 18
$ws.Range("R806").Value = "Hortaliza"

# New row 807
$ws.Range("A807").Value = 4
$ws.Range("B807").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C807").Value = "Los Lagos"
$ws.Range("D807").Value = 44918
$ws.Range("E807").Value = 10
$ws.Range("F807").Value = 100112020
$ws.Range("G807").Value = "Tomate"
$ws.Range("H807").Value = "Larga vida"
$ws.Range("I807").Value = "Primera"
$ws.Range("J807").Value = 250
$ws.Range("K807").Value = 22000
$ws.Range("L807").Value = 22000
$ws.Range("M807").Value = 22000
$ws.Range("N807").Value = "$/bandeja 18 kilos"
$ws.Range("O807").Value = "Provincia de Quillota"
$ws.Range("P807").Value = 1222
$ws.Range("Q807").Value = 18
$ws.Range("R807").Value = "Hortaliza"
